# Update "想去人数" (want-to-go count) values that changed when the site
# was regenerated (gh-pages output update at commit 456a3b4).
#
# Sheet "展览" (sheet1 / rId1):
#   F3:  417  -> 419
#   F12: 12089 -> 12091
#   F21: 359  -> 360
#
# Sheet "全部类型" (sheet4 / rId4) mirrors the same rows (shifted by +2 rows
# because it contains two extra rows before row 3's data starts repeating):
#   F3:  417  -> 419
#   F14: 12089 -> 12091
#   F23: 359  -> 360

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 419
$wsExhibit.Range("F12").Value = 12091
$wsExhibit.Range("F21").Value = 360

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 419
$wsAll.Range("F14").Value = 12091
$wsAll.Range("F23").Value = 360
